$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation is inserted above the existing two rows, so
# the current row 2 becomes row 3 and the current row 3 becomes row 4.
# Shift the existing data down one row at a time (bottom-up), cell by cell,
# using Value2 (plain read/write) and copying the date column's number
# format so the appearance is preserved without creating extra styles.
for ($col = 1; $col -le 18; $col++) {
    $ws.Cells.Item(4, $col).Value = $ws.Cells.Item(3, $col).Value2
}
$ws.Cells.Item(4, 4).NumberFormat = $ws.Cells.Item(3, 4).NumberFormat

for ($col = 1; $col -le 18; $col++) {
    $ws.Cells.Item(3, $col).Value = $ws.Cells.Item(2, $col).Value2
}
$ws.Cells.Item(3, 4).NumberFormat = $ws.Cells.Item(2, 4).NumberFormat

# New row 2: same market/category/etc, but the new week's Fecha and Volumen.
$ws.Range("A2").Value = 7
$ws.Range("B2").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C2").Value = "Ñuble"
$ws.Range("D2").Value = 44691
$ws.Range("E2").Value = 16
$ws.Range("F2").Value = 100112052
$ws.Range("G2").Value = "Albahaca"
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 3000
$ws.Range("L2").Value = 3500
$ws.Range("M2").Value = 3250
$ws.Range("N2").Value = "$/docena de matas"
$ws.Range("O2").Value = "Región Metropolitana"
$ws.Range("P2").Value = 542
$ws.Range("Q2").Value = 6
$ws.Range("R2").Value = "Hortaliza"
